# Update cryptos.xlsx per "Updated symbol list on Mon Dec 19 20:47:37 UTC 2022 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Price (column D) updates ---
Set-TextValue "D2"  "242.96"
Set-TextValue "D3"  "21.50"
Set-TextValue "D4"  "5.213"
Set-TextValue "D5"  "0.05602"
Set-TextValue "D7"  "6.375"
Set-TextValue "D8"  "0.8053"
Set-TextValue "D9"  "0.9527"
Set-TextValue "D10" "0.1433"
Set-TextValue "D11" "0.07305"
Set-TextValue "D12" "0.03142"
Set-TextValue "D14" "0.09280"
Set-TextValue "D15" "3.569"
Set-TextValue "D16" "0.001651"
Set-TextValue "D17" "0.04688"
Set-TextValue "D18" "0.0005749"
Set-TextValue "D19" "0.006360"
Set-TextValue "D20" "0.004985"
Set-TextValue "D22" "0.0001499"
Set-TextValue "D23" "0.0003100"
Set-TextValue "D24" "3.755"
Set-TextValue "D26" "0.3269"
Set-TextValue "D40" "0.03913"
Set-TextValue "D41" "0.006900"

# --- Row 42 & 43 swap (CEJI <-> BKEXToken) plus new prices ---
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1035"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002888"
$ws.Range("E43").Value = "42CEJICEJI"

# --- Remaining price (column D) updates ---
Set-TextValue "D44" "0.007506"
Set-TextValue "D45" "0.00005928"
Set-TextValue "D48" "0.6824"
Set-TextValue "D49" "0.07658"
